# Update the "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" worksheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 140
$ws1.Range("F5").Value  = 3320
$ws1.Range("F6").Value  = 1030
$ws1.Range("F7").Value  = 2191
$ws1.Range("F8").Value  = 2102
$ws1.Range("F9").Value  = 1101
$ws1.Range("F10").Value = 602
$ws1.Range("F13").Value = 393
$ws1.Range("F15").Value = 42
$ws1.Range("F18").Value = 1575
$ws1.Range("F19").Value = 626
$ws1.Range("F22").Value = 12210
$ws1.Range("F23").Value = 12244
$ws1.Range("F24").Value = 908
$ws1.Range("F28").Value = 18
$ws1.Range("F29").Value = 359
$ws1.Range("F31").Value = 193
$ws1.Range("F32").Value = 572

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 140
$ws4.Range("F6").Value  = 3320
$ws4.Range("F7").Value  = 1031
$ws4.Range("F8").Value  = 2191
$ws4.Range("F9").Value  = 2102
$ws4.Range("F10").Value = 1101
$ws4.Range("F11").Value = 602
$ws4.Range("F14").Value = 393
$ws4.Range("F17").Value = 42
$ws4.Range("F22").Value = 1575
$ws4.Range("F23").Value = 626
$ws4.Range("F26").Value = 12210
$ws4.Range("F27").Value = 12245
$ws4.Range("F28").Value = 908
$ws4.Range("F32").Value = 18
$ws4.Range("F33").Value = 359
$ws4.Range("F37").Value = 193
$ws4.Range("F38").Value = 572
